# Daily attendance processing - 2025-11-20 21:21:15
#
# The "Recorded By" column (G) contains comma-separated lists of names /
# email addresses (e.g. "dnasr281@gmail.com, System"). This script
# re-sorts each list using an ordinal (case-sensitive, byte-value) compare
# so the entries come out in a consistent, deterministic order.

function OrdinalLess($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minlen = $len1
    if ($len2 -lt $minlen) { $minlen = $len2 }
    for ($p_idx = 0; $p_idx -lt $minlen; $p_idx++) {
        $c1 = [int][char]$s1[$p_idx]
        $c2 = [int][char]$s2[$p_idx]
        if ($c1 -lt $c2) { return $true }
        if ($c1 -gt $c2) { return $false }
    }
    return ($len1 -lt $len2)
}

function SortOrdinal($arr) {
    $n = $arr.Count
    for ($sort_i = 1; $sort_i -lt $n; $sort_i++) {
        $key = $arr[$sort_i]
        $sort_j = $sort_i - 1
        while ($sort_j -ge 0 -and (OrdinalLess $key $arr[$sort_j])) {
            $arr[$sort_j + 1] = $arr[$sort_j]
            $sort_j = $sort_j - 1
        }
        $arr[$sort_j + 1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Value2
    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $names = $text -split ", "
    if ($names.Count -le 1) {
        continue
    }

    $sortedNames = SortOrdinal $names
    $newText = $sortedNames -join ", "
    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
